# Applies the "add latex reports and final results" edit:
#  - "Wyniki najlepszego" sheet: the three summary metric rows (correct /
#    f1_score / accuracy) are re-ordered so that f1_score, accuracy, correct
#    appear in that sequence (same label -> value associations, new order).
#  - Both cross-validation sheets ("Walidacja krzyzowa - trafnosc" and
#    "Walidacja krzyzowa - f1") get refreshed D-column (score) values for a
#    subset of folds/rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Wyniki najlepszego" (3rd sheet) - reorder the correct/f1_score/accuracy
#    rows (rows 3-5, column A = label, column B = value).
# ---------------------------------------------------------------------
$wsBest = $wb.Worksheets.Item(3)

$wsBest.Range("A3").Value = "f1_score"
$wsBest.Range("B3").Value = 0.8833094213295075

$wsBest.Range("A4").Value = "accuracy"
$wsBest.Range("B4").Value = 0.9019607843137255

$wsBest.Range("A5").Value = "correct"
$wsBest.Range("B5").Value = 0.9019607843137255

# ---------------------------------------------------------------------
# 2) Cross-validation results sheets - update column D values.
#    Sheets 4 ("Walidacja krzyzowa - trafnosc") and 5 ("Walidacja krzyzowa
#    - f1") receive identical updates to these rows.
# ---------------------------------------------------------------------
$dUpdates = @{
    4  = 0.8235294117647058
    5  = 0.7254901960784313
    6  = 0.7254901960784313
    7  = 0.7254901960784313
    8  = 0.7254901960784313
    10 = 0.7254901960784313
    11 = 0.7647058823529411
    12 = 0.7058823529411765
    13 = 0.7647058823529411
    18 = 0.6470588235294118
    24 = 0.6666666666666666
}

foreach ($sheetIndex in 4, 5) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $dUpdates.Keys) {
        $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
    }
}
